# "Cambios en primeros y siguientes"
# Update two cells in the Primeros/Siguientes table on Hoja1:
#  - D18 ("id" row): the "+,-" set now uses a single comma-joined string
#  - E23 ("id" row, Siguientes column): the Siguientes(expr_simple) set now
#    also includes "], oprel, ), ENTONCES, HACER, ;, SINO, "
# Also restore the selection left behind by the editing session (E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "id,num,(,NO,CIERTO,FALSO,+,-"
$ws.Range("E23").Value = "opsuma, O, ], oprel, ), ENTONCES, HACER, ;, SINO, "

$ws.Range("E7").Select() | Out-Null
